$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 18:35"

# Update Estados Unidos row (row 4)
$ws.Cells.Item(4, 2).Value = 1627786
$ws.Cells.Item(4, 3).Value = 6884
$ws.Cells.Item(4, 4).Value = 383678
$ws.Cells.Item(4, 5).Value = 1147399
$ws.Cells.Item(4, 7).Value = 355
$ws.Cells.Item(4, 8).Value = 96709

# Update Italia row (row 9)
$ws.Cells.Item(9, 2).Value = 228658
$ws.Cells.Item(9, 3).Value = 652
$ws.Cells.Item(9, 4).Value = 136720
$ws.Cells.Item(9, 5).Value = 59322
$ws.Cells.Item(9, 7).Value = 130
$ws.Cells.Item(9, 8).Value = 32616

# Update Argentina row (row 51)
$ws.Cells.Item(51, 4).Value = 3062
$ws.Cells.Item(51, 5).Value = 6450

# Update Maldivas row (row 101)
$ws.Cells.Item(101, 4).Value = 98
$ws.Cells.Item(101, 5).Value = 1114

# Rows 209-211 (Montserrat, Seychelles, Groenlandia) - the shared string
# table order changed (Groenlandia now precedes Montserrat & Seychelles)
# and the "casos activos" / "muertes" figures for rows 209/210 were updated.
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 0

$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Seychelles"
